# Lecture partielle de l'EDT M1 MIAGE.
# Shifts each week's date forward by 1096 days and updates the corresponding
# French weekday label to match the new date, and fixes a TP time typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new date serial, new weekday name)
$ws.Range("A2").Value2  = 46037
$ws.Range("B2").Value2  = "jeudi"

$ws.Range("A6").Value2  = 46038
$ws.Range("B6").Value2  = "vendredi"

$ws.Range("A10").Value2 = 46041
$ws.Range("B10").Value2 = "lundi"

$ws.Range("A13").Value2 = 46043
$ws.Range("B13").Value2 = "mercredi"

$ws.Range("A15").Value2 = 46045
$ws.Range("B15").Value2 = "vendredi"

$ws.Range("A19").Value2 = 46048
$ws.Range("B19").Value2 = "lundi"

$ws.Range("A22").Value2 = 46051
$ws.Range("B22").Value2 = "jeudi"

$ws.Range("A25").Value2 = 46052
$ws.Range("B25").Value2 = "vendredi"

$ws.Range("A28").Value2 = 46055
$ws.Range("B28").Value2 = "lundi"

$ws.Range("A31").Value2 = 46058
$ws.Range("B31").Value2 = "jeudi"

$ws.Range("A35").Value2 = 46080
$ws.Range("B35").Value2 = "vendredi"

$ws.Range("A38").Value2 = 46092
$ws.Range("B38").Value2 = "mercredi"

$ws.Range("A40").Value2 = 46094
$ws.Range("B40").Value2 = "vendredi"

$ws.Range("A44").Value2 = 46100
$ws.Range("B44").Value2 = "jeudi"

$ws.Range("A46").Value2 = 46101
$ws.Range("B46").Value2 = "vendredi"

$ws.Range("A50").Value2 = 46106
$ws.Range("B50").Value2 = "mercredi"

$ws.Range("A52").Value2 = 46107
$ws.Range("B52").Value2 = "jeudi"

$ws.Range("A54").Value2 = 46108
$ws.Range("B54").Value2 = "vendredi"

$ws.Range("A58").Value2 = 46113
$ws.Range("B58").Value2 = "mercredi"

$ws.Range("A60").Value2 = 46114
$ws.Range("B60").Value2 = "jeudi"

$ws.Range("A62").Value2 = 46115
$ws.Range("B62").Value2 = "vendredi"

# Fix TP start-time typo for the Monday (27/03) TP/GA session: 11:0 -> 10:0
$ws.Range("D53").Value2 = "10:0"
